$d = $word.ActiveDocument

$d.Content.Find.Execute("85×73=", $true, $false, $false, $false, $false, $true, 1, $false, "50×59=", 2) | Out-Null
$d.Content.Find.Execute("31×23=", $true, $false, $false, $false, $false, $true, 1, $false, "63×74=", 2) | Out-Null
$d.Content.Find.Execute("36×66=", $true, $false, $false, $false, $false, $true, 1, $false, "32×82=", 2) | Out-Null
$d.Content.Find.Execute("39×69=", $true, $false, $false, $false, $false, $true, 1, $false, "73×46=", 2) | Out-Null
$d.Content.Find.Execute("83×89=", $true, $false, $false, $false, $false, $true, 1, $false, "80×67=", 2) | Out-Null
$d.Content.Find.Execute("20×35=", $true, $false, $false, $false, $false, $true, 1, $false, "39×38=", 2) | Out-Null
$d.Content.Find.Execute("34×41=", $true, $false, $false, $false, $false, $true, 1, $false, "54×32=", 2) | Out-Null
$d.Content.Find.Execute("59×25=", $true, $false, $false, $false, $false, $true, 1, $false, "32×84=", 2) | Out-Null
$d.Content.Find.Execute("84×79=", $true, $false, $false, $false, $false, $true, 1, $false, "24×87=", 2) | Out-Null
$d.Content.Find.Execute("83×74=", $true, $false, $false, $false, $false, $true, 1, $false, "84×19=", 2) | Out-Null
$d.Content.Find.Execute("84×97=", $true, $false, $false, $false, $false, $true, 1, $false, "25×37=", 2) | Out-Null
$d.Content.Find.Execute("65×14=", $true, $false, $false, $false, $false, $true, 1, $false, "82×32=", 2) | Out-Null
$d.Content.Find.Execute("42×82=", $true, $false, $false, $false, $false, $true, 1, $false, "36×70=", 2) | Out-Null
$d.Content.Find.Execute("98×22=", $true, $false, $false, $false, $false, $true, 1, $false, "67×11=", 2) | Out-Null
$d.Content.Find.Execute("79×18=", $true, $false, $false, $false, $false, $true, 1, $false, "40×70=", 2) | Out-Null
$d.Content.Find.Execute("59×85=", $true, $false, $false, $false, $false, $true, 1, $false, "98×67=", 2) | Out-Null
$d.Content.Find.Execute("61×76=", $true, $false, $false, $false, $false, $true, 1, $false, "49×76=", 2) | Out-Null
$d.Content.Find.Execute("16×88=", $true, $false, $false, $false, $false, $true, 1, $false, "70×73=", 2) | Out-Null
$d.Content.Find.Execute("93×63=", $true, $false, $false, $false, $false, $true, 1, $false, "30×84=", 2) | Out-Null
$d.Content.Find.Execute("66×34=", $true, $false, $false, $false, $false, $true, 1, $false, "66×75=", 2) | Out-Null
$d.Content.Find.Execute("82×17=", $true, $false, $false, $false, $false, $true, 1, $false, "34×65=", 2) | Out-Null
$d.Content.Find.Execute("85×91=", $true, $false, $false, $false, $false, $true, 1, $false, "89×58=", 2) | Out-Null
$d.Content.Find.Execute("94×18=", $true, $false, $false, $false, $false, $true, 1, $false, "83×53=", 2) | Out-Null
$d.Content.Find.Execute("72×94=", $true, $false, $false, $false, $false, $true, 1, $false, "62×38=", 2) | Out-Null
$d.Content.Find.Execute("62×60=", $true, $false, $false, $false, $false, $true, 1, $false, "26×16=", 2) | Out-Null
